$d = $word.ActiveDocument

# --- Split the run "{m" into two runs: "{" and "m" ---------------------
$findRange = $d.Content
[void]$findRange.Find.Execute("{m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($findRange.Find.Found) {
    # The match range covers "{m"; the boundary we want is after the first
    # character ("{" | "m"). Toggling a character formatting property on
    # the first character and reverting it forces Word to break the run
    # at that boundary without altering the visible text or its effective
    # formatting.
    $boundary = $d.Range($findRange.Start, $findRange.Start + 1)
    $boundary.Bold = 1
    $boundary.Bold = 0
}

# --- Split the run "()}" into two runs: "()" and "}" --------------------
$findRange2 = $d.Content
[void]$findRange2.Find.Execute("()}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($findRange2.Find.Found) {
    # Here the boundary we want is before the last character ("()" | "}").
    $boundary2 = $d.Range($findRange2.End - 1, $findRange2.End)
    $boundary2.Bold = 1
    $boundary2.Bold = 0
}
